$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to Text format before writing values so that
# numeric-looking strings (e.g. "1.001", "0.3600") are preserved exactly as
# text and are not coerced into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.524.27'
$ws.Range("E2").Value = '  -2.88%  '

$ws.Range("D3").Value = '1.757.10'
$ws.Range("E3").Value = '  -3.46%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").Value = '323.75'
$ws.Range("E5").Value = '  -1.62%  '

$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.18%  '

$ws.Range("D7").Value = '0.4245'
$ws.Range("E7").Value = '  -3.03%  '

$ws.Range("D8").Value = '0.3600'
$ws.Range("E8").Value = '  -2.63%  '

$ws.Range("D9").Value = '0.07542'
$ws.Range("E9").Value = '  -2.59%  '

$ws.Range("D10").Value = '42.28'
$ws.Range("E10").Value = '  -6.20%  '

$ws.Range("D11").Value = '1.104'
$ws.Range("E11").Value = '  -3.45%  '

$ws.Range("E12").Value = '  +0.24%  '

$ws.Range("D13").Value = '20.78'
$ws.Range("E13").Value = '  -6.77%  '

$ws.Range("D14").Value = '6.059'
$ws.Range("E14").Value = '  -4.34%  '

$ws.Range("D15").Value = '7.200'
$ws.Range("E15").Value = '  -5.49%  '

$ws.Range("D16").Value = '1.754.60'
$ws.Range("E16").Value = '  -4.27%  '

$ws.Range("D17").Value = '92.87'
$ws.Range("E17").Value = '  -0.92%  '

$ws.Range("E18").Value = '  -1.56%  '

$ws.Range("D19").Value = '0.06385'
$ws.Range("E19").Value = '  -2.12%  '

$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.11%  '

$ws.Range("D21").Value = '17.06'
$ws.Range("E21").Value = '  -2.77%  '

$ws.Range("D22").Value = '5.904'
$ws.Range("E22").Value = '  -6.09%  '

$ws.Range("D23").Value = '27.572.63'
$ws.Range("E23").Value = '  -2.85%  '

$ws.Range("D24").Value = '11.27'
$ws.Range("E24").Value = '  -3.78%  '

$ws.Range("D25").Value = '2.115'
$ws.Range("E25").Value = '  -0.68%  '

$ws.Range("D26").Value = '162.19'
$ws.Range("E26").Value = '  +1.27%  '

$ws.Range("D27").Value = '20.22'
$ws.Range("E27").Value = '  -3.19%  '

$ws.Range("D28").Value = '1.955.48'
$ws.Range("E28").Value = '  -3.86%  '

$ws.Range("D29").Value = '2.165'
$ws.Range("E29").Value = '  -6.72%  '

$ws.Range("D30").Value = '125.45'
$ws.Range("E30").Value = '  -3.09%  '

$ws.Range("D31").Value = '1.102'
$ws.Range("E31").Value = '  -9.53%  '

$ws.Range("D32").Value = '5.597'
$ws.Range("E32").Value = '  -6.17%  '

$ws.Range("D33").Value = '3.651'
$ws.Range("E33").Value = '  +0.99%  '

$ws.Range("D34").Value = '0.08907'
$ws.Range("E34").Value = '  -3.45%  '

$ws.Range("D35").Value = '12.22'
$ws.Range("E35").Value = '  -5.75%  '

$ws.Range("E36").Value = '  -3.73%  '

$ws.Range("D37").Value = '0.2113'
$ws.Range("E37").Value = '  -3.63%  '

$ws.Range("D38").Value = '0.06019'
$ws.Range("E38").Value = '  -3.33%  '

$ws.Range("D39").Value = '0.6356'
$ws.Range("E39").Value = '  -4.03%  '

$ws.Range("D40").Value = '4.960'
$ws.Range("E40").Value = '  -4.66%  '

$ws.Range("D41").Value = '1.191'
$ws.Range("E41").Value = '  -0.37%  '

$ws.Range("D42").Value = '1.000'
$ws.Range("E42").Value = '  +0.11%  '

$ws.Range("D43").Value = '1.396'
$ws.Range("E43").Value = '  -2.14%  '

$ws.Range("D44").Value = '7.933'
$ws.Range("E44").Value = '  -3.49%  '

$ws.Range("D45").Value = '13.49'
$ws.Range("E45").Value = '  -3.16%  '

$ws.Range("D46").Value = '0.5922'
$ws.Range("E46").Value = '  -3.93%  '

$ws.Range("E47").Value = '  -1.40%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '1.987'
$ws.Range("E48").Value = '  -2.54%  '

$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '122.99'
$ws.Range("E49").Value = '  -3.02%  '

$ws.Range("D50").Value = '1.167'
$ws.Range("E50").Value = '  +0.54%  '

$ws.Range("D51").Value = '0.06842'
$ws.Range("E51").Value = '  -2.67%  '

# Restore normal (General) style for the Price column so no stray text
# formatting is left behind on cells after the values are written.
$ws.Range("D2:D51").Style = "Normal"